$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.561.39"
$ws.Range("E2").Value = "  -1.76%  "
$ws.Range("D3").Value = "2.062.96"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "242.53"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -7.84%  "
$ws.Range("D9").Value = "59.18"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("E10").Value = "  -7.08%  "
$ws.Range("D11").Value = "0.0752"
$ws.Range("E11").Value = "  -4.52%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("E14").Value = "  -9.38%  "
$ws.Range("D15").Value = "2.362.97"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "5.41"
$ws.Range("E16").Value = "  -6.07%  "
$ws.Range("D17").Value = "2.062.63"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "36.480.10"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").Value = "16.40"
$ws.Range("E19").Value = "  -11.98%  "
$ws.Range("D20").Value = "71.65"
$ws.Range("E20").Value = "  -4.30%  "
$ws.Range("D21").Value = "0.0₃0861"
$ws.Range("E21").Value = "  -4.45%  "
$ws.Range("D22").Value = "237.18"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").Value = "5.26"
$ws.Range("E23").Value = "  -4.32%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  -5.02%  "
$ws.Range("D26").Value = "9.44"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").Value = "2.13"
$ws.Range("E27").Value = "  -3.18%  "
$ws.Range("D28").Value = "164.28"
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("D29").Value = "20.42"
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("D31").Value = "5.06"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D33").Value = "4.58"
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("D34").Value = "0.0595"
$ws.Range("E34").Value = "  -4.86%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "2.27"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("E37").Value = "  +2.76%  "
$ws.Range("D38").Value = "0.0819"
$ws.Range("E38").Value = "  -8.01%  "
$ws.Range("E39").Value = "  -7.63%  "
$ws.Range("D40").Value = "2.92"
$ws.Range("E40").Value = "  -4.83%  "
$ws.Range("E41").Value = "  -8.15%  "
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("D43").Value = "0.0215"
$ws.Range("E43").Value = "  -3.98%  "
$ws.Range("D44").Value = "0.0934"
$ws.Range("E44").Value = "  -7.66%  "
$ws.Range("D45").Value = "94.27"
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("D46").Value = "1.386.34"
$ws.Range("E46").Value = "  +8.35%  "
$ws.Range("D47").Value = "7.44"
$ws.Range("E47").Value = "  +8.37%  "
$ws.Range("D48").Value = "15.53"
$ws.Range("E48").Value = "  -11.71%  "
$ws.Range("D49").Value = "2.35"
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("D50").Value = "2.85"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "2.251.88"
$ws.Range("E51").Value = "  +0.17%  "
